# Update "想去人数" (number of people interested) counts in column F
# for the rows whose values changed between the two data refreshes.
# The same data table exists (duplicated) on both the "展览" sheet and
# the "全部类型" sheet, so apply the updates to both.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 621
    3  = 2190
    4  = 78
    5  = 12942
    8  = 510
    9  = 474
    11 = 969
    12 = 13727
    13 = 14226
    22 = 1082
    26 = 5311
    28 = 288
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
